# "remove tokyo from all kpis"
#
# The "Functional KPIs" sheet carried an extra attribute pair (column O:
# store_attr_1_name, column P: store_attr_1_value) populated with
# "address_city" / "Tokyo" on every KPI row (rows 5-10). That data is being
# removed from all KPIs, leaving those cells blank (their existing
# formatting/style is kept).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional KPIs")

# Clear the "address_city" / "Tokyo" values from every KPI row.
$ws.Range("O5:P10").ClearContents()

# Leave the cursor where the author ended up after making the change.
$ws.Range("P19").Select() | Out-Null

# The workbook only has a single external reference, so the Validation_List
# defined name should point at external workbook index 1, not 2.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Validation_List") {
        $n.RefersTo = "=[1]Set_up!`$A`$90:`$A`$124"
    }
}
